# Update Ravarino.xlsx to 23 August 2021: append 14 new daily rows (344-357)
# covering dates 44418 (2021-08-10) .. 44431 (2021-08-23), each row holding
# [date serial, new positives, 7-day rolling sum, 7-day rolling sum per 100k].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows: row number, date serial (col A), col B, col C, col D
$newRows = @(
    @(344, 44418, 0, 3, 48.5201358563804),
    @(345, 44419, 0, 3, 48.5201358563804),
    @(346, 44420, 0, 1, 16.17337861879347),
    @(347, 44421, 0, 1, 16.17337861879347),
    @(348, 44422, 1, 1, 16.17337861879347),
    @(349, 44423, 2, 3, 48.5201358563804),
    @(350, 44424, 2, 5, 80.86689309396733),
    @(351, 44425, 0, 5, 80.86689309396733),
    @(352, 44426, 0, 5, 80.86689309396733),
    @(353, 44427, 0, 5, 80.86689309396733),
    @(354, 44428, 0, 5, 80.86689309396733),
    @(355, 44429, 1, 5, 80.86689309396733),
    @(356, 44430, 1, 4, 64.69351447517387),
    @(357, 44431, 0, 2, 32.34675723758694)
)

# The last existing row (343) carries the date cell style (centered, bordered,
# bold, date-formatted) that every subsequent date cell in column A re-uses.
$lastRow = 343

foreach ($r in $newRows) {
    $row = $r[0]

    # Clone column A's style from the row above so the new date cell matches
    # the rest of the column (same numFmt/border/font/alignment), then set
    # the real values on top of the pasted format.
    $ws.Range("A$lastRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Range("A$row").Value2 = $r[1]
    $ws.Range("B$row").Value2 = $r[2]
    $ws.Range("C$row").Value2 = $r[3]
    $ws.Range("D$row").Value2 = $r[4]

    $lastRow = $row
}

$excel.CutCopyMode = $false
